# "added corr direction to cookstove"
#
# Populate the new corr_direction column (L) on wrapper_ready for each of
# the four estimate rows, then autofit column A (its labels, e.g.
# "carbon_per_cookstove", no longer fully fit the default width) and leave
# the raw_data sheet's selection on the source rows (A2:A5) that back the
# new column.

$wb = $excel.ActiveWorkbook

$wsWrapper = $wb.Worksheets.Item("wrapper_ready")
$wsRaw = $wb.Worksheets.Item("raw_data")

$wsWrapper.Range("L2").Value = 1
$wsWrapper.Range("L3").Value = 1
$wsWrapper.Range("L4").Value = 2
$wsWrapper.Range("L5").Value = 2

# Widen column A to fit its (now more scrutinized) row labels.
$wsWrapper.Columns.Item(1).AutoFit() | Out-Null

# Highlight the raw_data rows (takeup_control, takeup_treatment, carbon
# reduction) used to determine each row's correlation direction.
$wsRaw.Activate()
$wsRaw.Range("A2:A5").Select() | Out-Null

# Leave the wrapper_ready sheet active/in front, as it was before editing,
# with the cursor back at the top (its prior H14 selection no longer
# reflects where the edit happened).
$wsWrapper.Activate()
$wsWrapper.Range("A1").Select() | Out-Null
